# Barriers_Pathway_Data.xlsx update
# - Removes the "Big Meadow Creek / Big Meadow Creek 01" reach row (row 3)
#   that is no longer part of the prioritization dataset.
# - Widens the ReachName column (C) now that it is the longest visible
#   column of text, and moves the active selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the entire "Big Meadow Creek" data row (Wenatchee basin), shifting
# the rows below it (Icicle Creek / Mission Creek entries) up by one.
$ws.Rows(3).Delete()

# Widen column C (ReachName) to fit the remaining entries.
$ws.Columns.Item(3).ColumnWidth = 28 + 1/6

# Move/restore the selection cursor.
$ws.Range("I11").Select()
